# Kobuleti Municipality area sheet - revert to the simpler export layout:
#  - drop the 1989 and 2002 census-year columns, keeping only the 2014 column
#  - drop the "(according to the population census data)" note row
#  - restore the larger 20.1pt row height used by this report layout
#  - keep four blank spacer rows below the data (matches the original export)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "1989" and "2002" data columns (B and C).
# The remaining "2014" column (old D) shifts left into column B,
# carrying its existing number formatting/borders with it.
$ws.Range("B:C").EntireColumn.Delete()

# Drop the whole row holding "(according to the population census data)".
# Everything below shifts up, so the former blank spacer row becomes row 2,
# "(sq. km)" becomes row 3, the year header becomes row 4, and the
# "Area" data row becomes row 5 - matching the target layout.
$ws.Range("A2").EntireRow.Delete()

# Restore the taller 20.1pt row height across the title/table rows plus the
# four trailing blank rows that round out the printed area.
$ws.Range("A1:B9").RowHeight = 20.1
